$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "First"
$ws.Range("D44").Value = "First"
$ws.Range("D45").Value = "First"
$ws.Range("D46").Value = "First"
$ws.Range("D47").Value = "First"
$ws.Range("D48").Value = "First"
$ws.Range("D49").Value = "Second"
$ws.Range("D50").Value = "Second"
$ws.Range("D51").Value = "Second"
$ws.Range("D52").Value = "First"
$ws.Range("D53").Value = "Second"
$ws.Range("D54").Value = "Second"
$ws.Range("D55").Value = "Second"
$ws.Range("D56").Value = "Second"
$ws.Range("D57").Value = "First"
$ws.Range("D59").Value = "Second"
$ws.Range("D60").Value = "Second"
$ws.Range("D61").Value = "First"
$ws.Range("D62").Value = "Second"
$ws.Range("D64").Value = "First"
$ws.Range("D65").Value = "Second"
$ws.Range("D67").Value = "First"
$ws.Range("D68").Value = "Second"
$ws.Range("D69").Value = "First"
$ws.Range("D70").Value = "First"
$ws.Range("D160").Value = "Bass"
$ws.Range("D161").Value = "Second"
$ws.Range("D162").Value = "Second"
$ws.Range("D163").Value = "First"
$ws.Range("D164").Value = "Second"
$ws.Range("D165").Value = "Bass"
$ws.Range("D166").Value = "First"
$ws.Range("D167").Value = "Second"
$ws.Range("D169").Value = "First"
$ws.Range("D170").Value = "Second"
$ws.Range("D171").Value = "First"
$ws.Range("D172").Value = "First"
$ws.Range("D173").Value = "Bass"
$ws.Range("D174").Value = "First"
$ws.Range("D175").Value = "First"
$ws.Range("D176").Value = "First"
$ws.Range("D177").Value = "First"
$ws.Range("D178").Value = "First"
$ws.Range("D179").Value = "Second"
$ws.Range("D180").Value = "Second"
$ws.Range("D181").Value = "Second"
$ws.Range("D182").Value = "First"
$ws.Range("D183").Value = "Second"
$ws.Range("D184").Value = "Second"
$ws.Range("D185").Value = "Second"
$ws.Range("D186").Value = "Second"
$ws.Range("D187").Value = "First"
$ws.Range("D189").Value = "Second"
$ws.Range("D191").Value = "Second"
$ws.Range("D193").Value = "Second"
$ws.Range("D194").Value = "Second"
$ws.Range("D196").Value = "Second"
$ws.Range("D198").Value = "Second"
$ws.Range("D201").Value = "Second"
$ws.Range("D202").Value = "First"
$ws.Range("D203").Value = "First"
$ws.Range("D204").Value = "Second"
$ws.Range("D205").Value = "Second"
$ws.Range("D206").Value = "First"
$ws.Range("D207").Value = "Second"
$ws.Range("D208").Value = "First"
$ws.Range("D209").Value = "First"
$ws.Range("D210").Value = "Second"
$ws.Range("D211").Value = "Second"
$ws.Range("D212").Value = "First"
$ws.Range("D213").Value = "Second"
$ws.Range("D214").Value = "Second"
$ws.Range("D215").Value = "First"
